$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, shifting existing rows 34:80 down to 35:81
$ws.Rows.Item(34).EntireRow.Insert()

# Populate the newly inserted row 34 with the new data record
$ws.Range("A34").Value = 1
$ws.Range("B34").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C34").Value = "Arica y Parinacota"
$ws.Range("D34").Value = 44658
$ws.Range("E34").Value = 15
$ws.Range("F34").Value = 100112038
$ws.Range("G34").Value = "Cebollín baby"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 300
$ws.Range("K34").Value = 1500
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = 1750
$ws.Range("N34").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O34").Value = "Región de Arica y Parinacota"
$ws.Range("P34").Value = 875
$ws.Range("Q34").Value = 2
$ws.Range("R34").Value = "Hortaliza"

# Ensure the date cell keeps the same date/time number format as the rest of column D
$ws.Range("D34").NumberFormat = $ws.Range("D35").NumberFormat
